$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 rows that were dropped from the catch table (delete bottom-up
# so earlier row numbers stay valid while iterating).
$ws.Rows(45).Delete()   # Wood NA
$ws.Rows(40).Delete()   # Shells NA
$ws.Rows(35).Delete()   # Liocarcinus depurator
$ws.Rows(32).Delete()   # Eggs of Murex
$ws.Rows(31).Delete()   # Biological discard

# All remaining catch weights (column G) for the 2-RAP gear rows are now
# reported as 0.
$ws.Range("G25:G40").Value = 0

$wb.Save()
